$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 408.3846
$ws.Range("I2").Value = 396.22223
$ws.Range("J2").Value = 435.75
$ws.Range("K2").Value = 396.22223
$ws.Range("L2").Value = 435.75
$ws.Range("M2").Value = -283.22223
$ws.Range("N2").Value = -661.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 353.45456
$ws.Range("I4").Value = 298.14285
$ws.Range("K4").Value = 298.14285
$ws.Range("M4").Value = -184.14285

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7491.533
$ws.Range("I43").Value = 7218.8
$ws.Range("J43").Value = 7627.9
$ws.Range("K43").Value = 7218.8
$ws.Range("L43").Value = 7627.9
$ws.Range("M43").Value = -7149.8
$ws.Range("N43").Value = -7765.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 192951.6
$ws.Range("I98").Value = 1140.5652
$ws.Range("K98").Value = 1140.5652
$ws.Range("M98").Value = 357.4348

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5391.7334
$ws.Range("I116").Value = 5041.5713
$ws.Range("J116").Value = 5698.125
$ws.Range("K116").Value = 5041.5713
$ws.Range("L116").Value = 5698.125
$ws.Range("M116").Value = -1599.5713
$ws.Range("N116").Value = -12582.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 810.8333
$ws.Range("J121").Value = 810.8333
$ws.Range("L121").Value = 2432.4999
$ws.Range("N121").Value = -5926.4999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 192951.6
$ws.Range("I122").Value = 1140.5652
$ws.Range("K122").Value = 3421.6956
$ws.Range("M122").Value = -971.6956

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1959.0834
$ws.Range("I132").Value = 1773.5454
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5320.6362
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -2790.6362
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2838.446
$ws.Range("I137").Value = 1441.1177
$ws.Range("J137").Value = 3333.3333
$ws.Range("K137").Value = 4323.3531
$ws.Range("L137").Value = 9999.999899999999
$ws.Range("M137").Value = -1773.3531
$ws.Range("N137").Value = -15099.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3107.6987
$ws.Range("I138").Value = 1999.6875
$ws.Range("J138").Value = 3418.7192
$ws.Range("K138").Value = 5999.0625
$ws.Range("L138").Value = 10256.1576
$ws.Range("M138").Value = -859.0625
$ws.Range("N138").Value = -20536.1576

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3339.85
$ws.Range("I32").Value = 2185.5535
$ws.Range("K32").Value = 2185.5535
$ws.Range("M32").Value = -1898.5535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2794.111
$ws.Range("I45").Value = 2794.111
$ws.Range("K45").Value = 2794.111
$ws.Range("M45").Value = -2417.111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6474.654
$ws.Range("I61").Value = 5782.2856
$ws.Range("J61").Value = 9382.6
$ws.Range("K61").Value = 5782.2856
$ws.Range("L61").Value = 9382.6
$ws.Range("M61").Value = -5570.2856
$ws.Range("N61").Value = -9806.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3022.75
$ws.Range("I132").Value = 2273.9614
$ws.Range("K132").Value = 6821.8842
$ws.Range("M132").Value = -4291.8842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6474.654
$ws.Range("I136").Value = 5782.2856
$ws.Range("J136").Value = 9382.6
$ws.Range("K136").Value = 17346.8568
$ws.Range("L136").Value = 28147.8
$ws.Range("M136").Value = -14796.8568
$ws.Range("N136").Value = -33247.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4155.893
$ws.Range("I20").Value = 3886.3157
$ws.Range("J20").Value = 4725
$ws.Range("K20").Value = 3886.3157
$ws.Range("L20").Value = 4725
$ws.Range("M20").Value = -3639.3157
$ws.Range("N20").Value = -5219

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3665
$ws.Range("I99").Value = 3495
$ws.Range("K99").Value = 3495
$ws.Range("M99").Value = -1997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26943.697
$ws.Range("I31").Value = 1641.3572
$ws.Range("J31").Value = 74174.734
$ws.Range("K31").Value = 1641.3572
$ws.Range("L31").Value = 74174.734
$ws.Range("M31").Value = -1346.3572
$ws.Range("N31").Value = -74764.734

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 26943.697
$ws.Range("I34").Value = 1641.3572
$ws.Range("J34").Value = 74174.734
$ws.Range("K34").Value = 1641.3572
$ws.Range("L34").Value = 74174.734
$ws.Range("M34").Value = -1439.3572
$ws.Range("N34").Value = -74578.734

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3368.7778
$ws.Range("J58").Value = 12504.667
$ws.Range("L58").Value = 12504.667
$ws.Range("N58").Value = -12910.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3324.8518
$ws.Range("I132").Value = 2630.24
$ws.Range("K132").Value = 7890.719999999999
$ws.Range("M132").Value = -5360.719999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2781.9565
$ws.Range("I134").Value = 1698.55
$ws.Range("K134").Value = 5095.65
$ws.Range("M134").Value = -2560.65

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3368.7778
$ws.Range("J136").Value = 12504.667
$ws.Range("L136").Value = 37514.001
$ws.Range("N136").Value = -42614.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3463.7646
$ws.Range("I34").Value = 2055
$ws.Range("J34").Value = 4449.9
$ws.Range("K34").Value = 6165
$ws.Range("L34").Value = 13349.7
$ws.Range("M34").Value = -6081
$ws.Range("N34").Value = -13517.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1333.7142
$ws.Range("I113").Value = 1095.5
$ws.Range("J113").Value = 1550.2727
$ws.Range("K113").Value = 3286.5
$ws.Range("L113").Value = 4650.8181
$ws.Range("M113").Value = -1116.5
$ws.Range("N113").Value = -8990.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3953.625
$ws.Range("I139").Value = 2866.3635
$ws.Range("J139").Value = 6345.6
$ws.Range("K139").Value = 8599.0905
$ws.Range("L139").Value = 19036.8
$ws.Range("M139").Value = -3459.0905
$ws.Range("N139").Value = -29316.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14145.308
$ws.Range("I70").Value = 13413.286
$ws.Range("K70").Value = 13413.286
$ws.Range("M70").Value = -13143.286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 14145.308
$ws.Range("I73").Value = 13413.286
$ws.Range("K73").Value = 13413.286
$ws.Range("M73").Value = -12477.286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 181840.83
$ws.Range("I80").Value = 280027
$ws.Range("J80").Value = 5105.7
$ws.Range("K80").Value = 280027
$ws.Range("L80").Value = 5105.7
$ws.Range("M80").Value = -279029
$ws.Range("N80").Value = -7101.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 181840.83
$ws.Range("I83").Value = 280027
$ws.Range("J83").Value = 5105.7
$ws.Range("K83").Value = 1400135
$ws.Range("L83").Value = 25528.5
$ws.Range("M83").Value = -1395143
$ws.Range("N83").Value = -35512.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2191.3928
$ws.Range("I97").Value = 1297.8695
$ws.Range("K97").Value = 1297.8695
$ws.Range("M97").Value = -801.8695

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8505.657999999999
$ws.Range("I122").Value = 9510.727999999999
$ws.Range("K122").Value = 28532.184
$ws.Range("M122").Value = -26082.184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2898.1035
$ws.Range("I132").Value = 2100.625
$ws.Range("K132").Value = 6301.875
$ws.Range("M132").Value = -3771.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5789.811
$ws.Range("I7").Value = 4244.1724
$ws.Range("J7").Value = 11392.75
$ws.Range("K7").Value = 4244.1724
$ws.Range("L7").Value = 11392.75
$ws.Range("M7").Value = -4132.1724
$ws.Range("N7").Value = -11616.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4334
$ws.Range("I16").Value = 1001
$ws.Range("J16").Value = 6000.5
$ws.Range("K16").Value = 1001
$ws.Range("L16").Value = 6000.5
$ws.Range("M16").Value = -831
$ws.Range("N16").Value = -6340.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7835.933
$ws.Range("I40").Value = 6723.85
$ws.Range("K40").Value = 6723.85
$ws.Range("M40").Value = -6587.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5302.676
$ws.Range("I46").Value = 3793.6875
$ws.Range("J46").Value = 6452.381
$ws.Range("K46").Value = 3793.6875
$ws.Range("L46").Value = 6452.381
$ws.Range("M46").Value = -3605.6875
$ws.Range("N46").Value = -6828.381

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 500
$ws.Range("I50").Value = 500
$ws.Range("K50").Value = 500
$ws.Range("M50").Value = 137

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2099.2812
$ws.Range("I55").Value = 1243.2778
$ws.Range("J55").Value = 3199.8572
$ws.Range("K55").Value = 1243.2778
$ws.Range("L55").Value = 3199.8572
$ws.Range("M55").Value = -1070.2778
$ws.Range("N55").Value = -3545.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 44333.332
$ws.Range("J110").Value = 44333.332
$ws.Range("L110").Value = 44333.332
$ws.Range("N110").Value = -52513.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5789.811
$ws.Range("I126").Value = 4244.1724
$ws.Range("J126").Value = 11392.75
$ws.Range("K126").Value = 12732.5172
$ws.Range("L126").Value = 34178.25
$ws.Range("M126").Value = -10262.5172
$ws.Range("N126").Value = -39118.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5191.0435
$ws.Range("I132").Value = 3937.1538
$ws.Range("J132").Value = 6821.1
$ws.Range("K132").Value = 11811.4614
$ws.Range("L132").Value = 20463.3
$ws.Range("M132").Value = -9281.4614
$ws.Range("N132").Value = -25523.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 10000
$ws.Range("K25").Value = 10000
$ws.Range("M25").Value = -9707

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5285.1113
$ws.Range("I81").Value = 3446.4
$ws.Range("J81").Value = 5992.3076
$ws.Range("K81").Value = 6892.8
$ws.Range("L81").Value = 11984.6152
$ws.Range("M81").Value = -5831.8
$ws.Range("N81").Value = -14106.6152

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5285.1113
$ws.Range("I84").Value = 3446.4
$ws.Range("J84").Value = 5992.3076
$ws.Range("K84").Value = 34464
$ws.Range("L84").Value = 59923.076
$ws.Range("M84").Value = -29160
$ws.Range("N84").Value = -70531.076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 50073.75
$ws.Range("J119").Value = 50073.75
$ws.Range("L119").Value = 50073.75
$ws.Range("N119").Value = -59749.75
